# Added test case scenario to verify user type
# Update the SmokeTest flag ("D" column) for the first two test cases
# (TC_001, TC_002) from "No" to "Yes".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Yes"
$ws.Range("D3").Value = "Yes"

# Move the active selection to D3, matching the cursor position left
# behind after making this edit.
$ws.Range("D3").Select()
